$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q4" right after "总计", by duplicating
#    the "2022-Q3" sheet (so headers/styles/column widths all match the
#    existing per-quarter sheets), then overwriting its data with the
#    new quarter's single fund row.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)
$wsQ3    = $wb.Worksheets.Item("2022-Q3")

$wsQ3.Copy($null, $wsTotal)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# The copied sheet has 11 data rows (rows 2-12); the new quarter only has
# one, so clear rows 3-12 entirely (values + formatting).
$newSheet.Range("A3:H12").Clear()

# Fill in the single data row for 2022-Q4.
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "'012495"
$newSheet.Cells.Item(2,3).Value = "'民生加银双核动力混合"
$newSheet.Cells.Item(2,4).Value = "'0.12"
$newSheet.Cells.Item(2,5).Value = "'93.75"
$newSheet.Cells.Item(2,6).Value = "'4.88"
$newSheet.Cells.Item(2,7).Value = "'0.0059"
$newSheet.Cells.Item(2,8).Value = 6

# Writing text via a leading apostrophe forces a text-quote-prefix style;
# strip that back to the plain (unstyled) look used elsewhere by pasting
# the format from an already-cleared, default-styled cell.
$newSheet.Cells.Item(3,2).Copy()
$newSheet.Range("B2:G2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Insert the corresponding summary row into "总计" (sheet 1), pushing
#    the existing quarters down by one row and relabelling/renumbering
#    them accordingly.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(7,1).Value = 5
$ws1.Cells.Item(6,1).Copy()
$ws1.Cells.Item(7,1).PasteSpecial(-4122)
$ws1.Cells.Item(7,2).Value = "2021-Q3"
$ws1.Cells.Item(7,3).Value = 1
$ws1.Cells.Item(7,4).Value = 2.44

$ws1.Cells.Item(6,1).Value = 4
$ws1.Cells.Item(6,2).Value = "2021-Q4"
$ws1.Cells.Item(6,3).Value = 3
$ws1.Cells.Item(6,4).Value = 4.33

$ws1.Cells.Item(5,1).Value = 3
$ws1.Cells.Item(5,2).Value = "2022-Q1"
$ws1.Cells.Item(5,3).Value = 7
$ws1.Cells.Item(5,4).Value = 5.21

$ws1.Cells.Item(4,1).Value = 2
$ws1.Cells.Item(4,2).Value = "2022-Q2"
$ws1.Cells.Item(4,3).Value = 6
$ws1.Cells.Item(4,4).Value = 5.33

$ws1.Cells.Item(3,1).Value = 1
$ws1.Cells.Item(3,2).Value = "2022-Q3"
$ws1.Cells.Item(3,3).Value = 11
$ws1.Cells.Item(3,4).Value = 0.47

$ws1.Cells.Item(2,1).Value = 0
$ws1.Cells.Item(2,2).Value = "2022-Q4"
$ws1.Cells.Item(2,3).Value = 1
$ws1.Cells.Item(2,4).Value = 0.01
